$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Cells.Item(200,200)

$ws.Cells.Item(2,1).Value = 'Pokemon'
$ws.Cells.Item(2,2).Value = 'Pokemon.com: The Official Pokémon Website'
$ws.Cells.Item(2,3).Value = 'https://www.pokemon.com/us'
$ws.Cells.Item(2,4).Value = 'What''s New This Week · Dragon-Type Pokémon Soar over Pokéween · Scarlet & Violet—Stellar Crown Triple Play Deck Strategy · Live, Laugh, Scare with Pokémon ...'

$ws.Cells.Item(3,1).Value = 'Pokemon'
$ws.Cells.Item(3,2).Value = 'Pokémon'
$ws.Cells.Item(3,3).Value = 'https://en.wikipedia.org/wiki/Pok%C3%A9mon'
$ws.Cells.Item(3,4).Value = 'Pokémon is a Japanese media franchise consisting of video games, animated series and films, a trading card game, and other related media.'

$ws.Cells.Item(4,1).Value = 'Pokemon'
$ws.Cells.Item(4,2).Value = 'Pokémon GO'
$ws.Cells.Item(4,3).Value = 'https://pokemongolive.com/'
$ws.Cells.Item(4,4).Value = 'Spin Photo Discs at PokéStops and Gyms to receive items that will help you on your Pokémon GO journey, including Poké Balls, Berries, Evolution items, and more.'

$ws.Cells.Item(5,1).Value = 'Pokemon'
$ws.Cells.Item(5,2).Value = 'The Official Pokémon YouTube channel'
$ws.Cells.Item(5,3).Value = 'https://www.youtube.com/channel/UCFctpiB_Hnlk3ejWfHqSm6Q'
$ws.Cells.Item(5,4).Value = 'Watch Pokémon Video Game, Trading Card Game, and Animation content, along with special features, and event highlights right here on the official Pokémon ...'

$ws.Cells.Item(6,1).Value = 'Pokemon'
$ws.Cells.Item(6,2).Value = 'Pokemon - Pokémon'
$ws.Cells.Item(6,3).Value = 'https://twitter.com/pokemon'
$ws.Cells.Item(6,4).Value = 'Infernape with the Mightiest Mark is now appearing in black crystal Tera Raid Battles throughout #PokemonScarletViolet! Work together with friends to topple ...'

$ws.Cells.Item(7,1).Value = 'Minecraft'
$ws.Cells.Item(7,2).Value = 'Welcome to the Minecraft Official Site | Minecraft'
$ws.Cells.Item(7,3).Value = 'https://www.minecraft.net/'
$ws.Cells.Item(7,4).Value = 'Explore new gaming adventures, accessories, & merchandise on the Minecraft Official Site. Buy & download the game here, or check the site for the latest ...'

$ws.Cells.Item(8,1).Value = 'Minecraft'
$ws.Cells.Item(8,2).Value = 'Minecraft'
$ws.Cells.Item(8,3).Value = 'https://en.wikipedia.org/wiki/Minecraft'
$ws.Cells.Item(8,4).Value = 'Minecraft is a 2011 sandbox game developed and published by Swedish video game developer Mojang Studios. Originally created by Markus "Notch" Persson using ...'

$ws.Cells.Item(9,1).Value = 'Minecraft'
$ws.Cells.Item(9,2).Value = 'Minecraft: Play with Friends - Apps on Google Play'
$ws.Cells.Item(9,3).Value = 'https://play.google.com/store/apps/details?id=com.mojang.minecraftpe&hl=en_US'
$ws.Cells.Item(9,4).Value = '— Explore and craft your way through a completely open world where you can play with friends, build a city, start a farm, mine deep into the ...'

$ws.Cells.Item(10,1).Value = 'Minecraft'
$ws.Cells.Item(10,2).Value = 'Minecraft - PS4 Games'
$ws.Cells.Item(10,3).Value = 'https://www.playstation.com/en-ie/games/minecraft/'
$ws.Cells.Item(10,4).Value = 'Team up or go solo and triumph over waves of hostile mobs, build with new blocks, harness auto-crafting, battle the breeze, unlock the vault, and more.'

$ws.Cells.Item(11,1).Value = 'Minecraft'
$ws.Cells.Item(11,2).Value = 'r/Minecraft'
$ws.Cells.Item(11,3).Value = 'https://www.reddit.com/r/Minecraft/'
$ws.Cells.Item(11,4).Value = 'r/Minecraft: Minecraft community on Reddit.'

$ws.Cells.Item(12,1).Value = 'Burguer'
$ws.Cells.Item(12,2).Value = 'Burger King'
$ws.Cells.Item(12,3).Value = 'https://www.burgerking.pt/pt/'
$ws.Cells.Item(12,4).Value = 'Nearby restaurant! Free delivery. Free Delivery. In orders up to 20€. Order now!'

$ws.Cells.Item(13,1).Value = 'Burguer'
$ws.Cells.Item(13,2).Value = 'Burger King®'
$ws.Cells.Item(13,3).Value = 'https://www.burgerking.com.br/'
$ws.Cells.Item(13,4).Value = 'Baixe nosso App e tenha o BK na palma da sua mão! · Centro de preferências de privacidade.'

$ws.Cells.Item(14,1).Value = 'Burguer'
$ws.Cells.Item(14,2).Value = 'Home - Burger KingBurger King | HUNGRY? WE GOT YOU'
$ws.Cells.Item(14,3).Value = 'https://www.whopper.ie/'
$ws.Cells.Item(14,4).Value = 'Gourmet Kings · Veggie & Plant-based Kings · Texas Bacon Lovers · Sweet Treats.'

$ws.Cells.Item(15,1).Value = 'Burguer'
$ws.Cells.Item(15,2).Value = 'Bilbo Burguer ®️ (@bilboburguer)'
$ws.Cells.Item(15,3).Value = 'https://www.instagram.com/bilboburguer/'
$ws.Cells.Item(15,4).Value = 'CUPOM PRIMEIRA COMPRA: PRIMEIRAVEZ 🎟️ R. José Moreira, Itinga em Lauro de Freitas Pedindo no site chega mais rápido ⇊ Link para pedidos abaixo.'

$ws.Cells.Item(16,1).Value = 'Burguer'
$scratch.Formula = '="''O'' Burguer"'
$scratch.Copy()
$ws.Cells.Item(16,2).PasteSpecial(-4163)
$ws.Cells.Item(16,3).Value = 'http://www.oburguer.com.br/'
$scratch.Formula = '="''O'' Burguer é muito mais do que uma simples hamburgueria. Inspirada pela ciência e pela criatividade gastronômica, cada um de nossos produtos é uma experiência ..."'
$scratch.Copy()
$ws.Cells.Item(16,4).PasteSpecial(-4163)

$ws.Cells.Item(17,1).Value = 'League of Legends'
$ws.Cells.Item(17,2).Value = 'League of Legends Homepage'
$ws.Cells.Item(17,3).Value = 'https://www.leagueoflegends.com/'
$ws.Cells.Item(17,4).Value = 'League of Legends is a team-based game with over 140 champions to make epic plays with. Play now for free.'

$ws.Cells.Item(18,1).Value = 'League of Legends'
$ws.Cells.Item(18,2).Value = 'League of Legends'
$ws.Cells.Item(18,3).Value = 'https://en.wikipedia.org/wiki/League_of_Legends'
$ws.Cells.Item(18,4).Value = 'League of Legends (LoL), commonly referred to as League, is a 2009 multiplayer online battle arena video game developed and published by Riot Games.'

$ws.Cells.Item(19,1).Value = 'League of Legends'
$ws.Cells.Item(19,2).Value = 'LoL Esports | SCHEDULE'
$ws.Cells.Item(19,3).Value = 'https://lolesports.com/'
$ws.Cells.Item(19,4).Value = 'The best place to watch LoL Esports and earn rewards!'

$ws.Cells.Item(20,1).Value = 'League of Legends'
$ws.Cells.Item(20,2).Value = 'League of Legends'
$ws.Cells.Item(20,3).Value = 'https://www.youtube.com/channel/UC2t5bjwHdUX4vM2g8TRDq5g'
$ws.Cells.Item(20,4).Value = 'the most played video game in the world—100 million play every month—League of Legends® is a multiplayer online battle arena game (MOBA).'

$ws.Cells.Item(21,1).Value = 'League of Legends'
$ws.Cells.Item(21,2).Value = 'Riot Games. Developer of League of Legends, VALORANT ...'
$ws.Cells.Item(21,3).Value = 'https://www.riotgames.com/'
$ws.Cells.Item(21,4).Value = 'Riot Games. Developer of League of Legends, VALORANT, Teamfight Tactics, Legends of Runeterra, and Wild Rift. Creators of Arcane. Home of LOL and VALORANT ...'

$ws.Cells.Item(22,1).Value = 'Wild Rift'
$ws.Cells.Item(22,2).Value = 'Wild Rift: Welcome to League of Legends'
$ws.Cells.Item(22,3).Value = 'https://wildrift.leagueoflegends.com/'
$ws.Cells.Item(22,4).Value = 'League of Legends: Wild Rift - Team up with friends and test your skills in 5v5 mobile MOBA combat.'

$ws.Cells.Item(23,1).Value = 'Wild Rift'
$ws.Cells.Item(23,2).Value = 'League of Legends: Wild Rift - Apps on Google Play'
$ws.Cells.Item(23,3).Value = 'https://play.google.com/store/apps/details?id=com.riotgames.league.wildrift'
$ws.Cells.Item(23,4).Value = '— Enjoy fast-paced MOBA combat, real-time strategy, smooth controls, and diverse 5v5 gameplay. Team up with friends, lock in your champion, and ...'

$ws.Cells.Item(24,1).Value = 'Wild Rift'
$ws.Cells.Item(24,2).Value = 'League of Legends: Wild Rift'
$ws.Cells.Item(24,3).Value = 'https://www.youtube.com/c/wildrift'
$ws.Cells.Item(24,4).Value = 'Dive into League of Legends: Wild Rift: the skills-and-strategy 5v5 MOBA experience of League of Legends by Riot Games, now built from the ground up for ...'

$ws.Cells.Item(25,1).Value = 'Wild Rift'
$ws.Cells.Item(25,2).Value = 'WildRiftFire: Wild Rift Builds & Guides'
$ws.Cells.Item(25,3).Value = 'https://www.wildriftfire.com/'
$ws.Cells.Item(25,4).Value = 'Find the best build guide for Wild Rift on WildRiftFire. Learn which items, runes, and summoner spells to take on each champion in Wild Rift, as well as how ...'

$ws.Cells.Item(26,1).Value = 'Wild Rift'
$ws.Cells.Item(26,2).Value = 'League of Legends: Wild Rift (@wildrift) ...'
$ws.Cells.Item(26,3).Value = 'https://twitter.com/wildrift'
$ws.Cells.Item(26,4).Value = 'Official account of League of Legends: Wild Rift, your favorite mobile MOBA from @riotgames.'

$scratch.ClearContents()
$excel.CutCopyMode = $false